# Mejorar formato de trabajos de grado dirigidos:
# add a "why" (thesis) column entry for the two doctoral/master's rows,
# with wrapped text, taller rows, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tesisFarc = "Tesis: \href{http://hdl.handle.net/10234/167275}{Implementar una gu" + [char]0x00ED + "a para el seguimiento period" + [char]0x00ED + "stico con enfoque de g" + [char]0x00E9 + "nero de las desmovilizaciones de las FARC}"
$tesisCumLaude = "Tesis: \href{http://hdl.handle.net/10016/35862}{El tratamiento period" + [char]0x00ED + "stico de la violencia sexual en contra de las mujeres en el marco del conflicto armado colombiano. An" + [char]0x00E1 + "lisis de casos seg" + [char]0x00FA + "n el tipo de violencia, v" + [char]0x00ED + "ctimas, victimarios y contextos}. \textit{Cum Laude} y menci" + [char]0x00F3 + "n Internacional"

# Register the FARC thesis text first so it takes the lower shared-string
# index, then the Cum Laude text, matching the source workbook's string table
# order. Row 3 (Maestria, 2017) gets the FARC text; Row 2 (Doctorado, 2022)
# gets the Cum Laude text.
$ws.Range("E3").Value = $tesisFarc
$ws.Range("E2").Value = $tesisCumLaude

# Wrap text for the new cells and make the rows taller
$ws.Range("E2:E3").WrapText = $true
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30

# Move the active selection to E2
$ws.Range("E2").Select()

$wb.Save()
